# Actualización automática 2025-09-30 15:30:09
# Applies updated September ("septiembre") sales figures across the three
# report sheets: VENTAS POR GRUPO (per-group breakdown), VENTA MENSUAL
# (per-client monthly totals) and CUMPLIMIENTO MENSUAL (budget-compliance
# summary), plus the "N de 53" coverage-count row on VENTAS POR GRUPO.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M4").Value = 2521.53
$wsGrupo.Range("M5").Value = 3887.15
$wsGrupo.Range("M21").Value = 1054.31
$wsGrupo.Range("L22").Value = 2780.75
$wsGrupo.Range("M22").Value = 4359.41
$wsGrupo.Range("H27").Value = 1161
$wsGrupo.Range("I27").Value = 43.2
$wsGrupo.Range("I34").Value = 626.4
$wsGrupo.Range("P34").Value = 316.28
$wsGrupo.Range("L46").Value = 1140.48
$wsGrupo.Range("M46").Value = 3999.25
$wsGrupo.Range("C47").Value = 518.4

# Row 55 holds "<n> de 53" coverage counters per column; five of them
# increment by one because the cells above newly went from 0 to non-zero.
$wsGrupo.Range("C55").Value = "3 de 53"
$wsGrupo.Range("I55").Value = "11 de 53"
$wsGrupo.Range("L55").Value = "10 de 53"
$wsGrupo.Range("M55").Value = "18 de 53"
$wsGrupo.Range("P55").Value = "2 de 53"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL (column F = septiembre)
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value = 6849.29
$wsMensual.Range("F5").Value = 6679.65
$wsMensual.Range("F21").Value = 1054.31
$wsMensual.Range("F22").Value = 8432.22
$wsMensual.Range("F27").Value = 3417.52
$wsMensual.Range("F34").Value = 11589.38
$wsMensual.Range("F46").Value = 6441.13
$wsMensual.Range("F47").Value = 518.4
$wsMensual.Range("F59").Value = 86871.93

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL (C=PRESUPUESTO, D=VENTA, E=POR CUMPLIR,
# F=CUMPLIMIENTO). Only D/E/F move; C (budget) is unchanged.
# ---------------------------------------------------------------------
$wsCumple = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumple.Range("D2").Value = 1321.92
$wsCumple.Range("E2").Value = 4875.66402943659
$wsCumple.Range("F2").Value = 0.2132960188552979

$wsCumple.Range("D6").Value = 4072.5
$wsCumple.Range("E6").Value = -1164.91631853974
$wsCumple.Range("F6").Value = 1.400647563806209

$wsCumple.Range("D7").Value = 2753.1
$wsCumple.Range("E7").Value = -1866.388983712426
$wsCumple.Range("F7").Value = 3.104844700730691

$wsCumple.Range("D8").Value = 316.28
$wsCumple.Range("E8").Value = 1030.12488751609
$wsCumple.Range("F8").Value = 0.2349070498276993

$wsCumple.Range("D11").Value = 15667.05
$wsCumple.Range("E11").Value = 2164.364398465401
$wsCumple.Range("F11").Value = 0.878620711172992

$wsCumple.Range("D12").Value = 51474.1
$wsCumple.Range("E12").Value = 10389.6203947566
$wsCumple.Range("F12").Value = 0.8320563275461008

$wsCumple.Range("D15").Value = 84509.07
$wsCumple.Range("E15").Value = 37545.76551083435
$wsCumple.Range("F15").Value = 0.6923860873377561
